$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Eventos")

$ws.Range("K2").Value = 'https://s3.samsara.com/samsara-dashcam-videos/4006124/281474991205821/1748891365835/F7KsVXPM8Y-camera-video-segment-driver-1748891368335.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSB6HHK3ZB%2F20250603%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250603T190046Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEEMaCXVzLXdlc3QtMiJHMEUCIQDKu6SVHCuPxyLrQICvTWjsc8bJ1CoUDXThxrO8uE40HgIgA430WHEcncckntaqxq6OjxUxQq%2F4cRp0Tc2i6CEzSzcq3QMIGxAEGgw3ODEyMDQ5NDIyNDQiDI9ylrE874qsoNZPVyq6AyjSTgVQ3WieFhqFryrKxkcF3ppOLD9zcpUrjVHb3tfkBN6XzMndrmHyS91MRE81cdbF8AgqQTfi0m2KmkyBJ6n4LzAczOL5eQmIGQqjRTfVfv4MbIJNcX%2FmVKFDY46VPtANd2dKZKWtvOAG2jurMALhCaXaFkW38Wi%2FDpDX4kGeQYxsZ%2BWqA%2FsIP1TUnMfrtTnXOXr2FGA%2FMtA8a3pU%2B9b5iTKTFU1q8iQeuEehHKOXDFejJzr6rd39mMAWUoBhL7ZoAfvRYwq%2FK%2Bu5qH2XGR4aTPF6Dx2xIKkAwnhQPUa0qE2CeRnfYOMJhWzPQZ%2FCAldov9ZAih%2FeGjhBkztd%2BwikotUi51k1gb5djBDdD2%2BnrfS7%2FA%2BTSmkoMhyATJSRIqetQy1V%2B9I62PDv6sKYPAk1LKJnxA9P2VzZjctcIoYlmqoLPRDc2V8iEP0%2FlBt2b6DkkgBhIYGN%2Bc5rO4G%2FDic5glV%2BC219thdfzs8tXE%2BWfqMdD%2BC6aDJnVkR5UhAxhc7sAU4AkjlsvBnyBmxe132qO3Z%2FiJXr5kbp5dblgunXwgVgK9j%2B8yggr%2BiTKBbxM48Jy6PHKAR5CS4wyvr8wQY6pQGka52Q1n8aqPwGgfIPuwP%2B0wrpKmLdlRxGHPY4UeSOY6kBGwsVeZJGugNXCjDZUPIt8Kbou5blV%2Fh4c21o97PMLU%2FWg4Uln96yBi32j2p9lYTYOYzY6hmHQZ3%2Br1IP6tipUXxLrSP0lJlUGMnfhDJrO4lrGDrfWRWkkBBLnKeb%2FXTdLeeqsmurU7%2Foh%2FG2kJZoNWJTiGsVPbnh48Wng6XwbEgLefI%3D&X-Amz-SignedHeaders=host&response-expires=Wed%2C%2004%20Jun%202025%2003%3A00%3A46%20GMT&X-Amz-Signature=1fb23d771edee7338a15794ef1df1b0e5a25ca68c57b7356b7177b729d7b5b88'

$ws.Range("L2").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991205821/1748891365835/MSujKlorgb-camera-video-segment-1748891368335.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSB6HHK3ZB%2F20250603%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250603T190046Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEEMaCXVzLXdlc3QtMiJHMEUCIQDKu6SVHCuPxyLrQICvTWjsc8bJ1CoUDXThxrO8uE40HgIgA430WHEcncckntaqxq6OjxUxQq%2F4cRp0Tc2i6CEzSzcq3QMIGxAEGgw3ODEyMDQ5NDIyNDQiDI9ylrE874qsoNZPVyq6AyjSTgVQ3WieFhqFryrKxkcF3ppOLD9zcpUrjVHb3tfkBN6XzMndrmHyS91MRE81cdbF8AgqQTfi0m2KmkyBJ6n4LzAczOL5eQmIGQqjRTfVfv4MbIJNcX%2FmVKFDY46VPtANd2dKZKWtvOAG2jurMALhCaXaFkW38Wi%2FDpDX4kGeQYxsZ%2BWqA%2FsIP1TUnMfrtTnXOXr2FGA%2FMtA8a3pU%2B9b5iTKTFU1q8iQeuEehHKOXDFejJzr6rd39mMAWUoBhL7ZoAfvRYwq%2FK%2Bu5qH2XGR4aTPF6Dx2xIKkAwnhQPUa0qE2CeRnfYOMJhWzPQZ%2FCAldov9ZAih%2FeGjhBkztd%2BwikotUi51k1gb5djBDdD2%2BnrfS7%2FA%2BTSmkoMhyATJSRIqetQy1V%2B9I62PDv6sKYPAk1LKJnxA9P2VzZjctcIoYlmqoLPRDc2V8iEP0%2FlBt2b6DkkgBhIYGN%2Bc5rO4G%2FDic5glV%2BC219thdfzs8tXE%2BWfqMdD%2BC6aDJnVkR5UhAxhc7sAU4AkjlsvBnyBmxe132qO3Z%2FiJXr5kbp5dblgunXwgVgK9j%2B8yggr%2BiTKBbxM48Jy6PHKAR5CS4wyvr8wQY6pQGka52Q1n8aqPwGgfIPuwP%2B0wrpKmLdlRxGHPY4UeSOY6kBGwsVeZJGugNXCjDZUPIt8Kbou5blV%2Fh4c21o97PMLU%2FWg4Uln96yBi32j2p9lYTYOYzY6hmHQZ3%2Br1IP6tipUXxLrSP0lJlUGMnfhDJrO4lrGDrfWRWkkBBLnKeb%2FXTdLeeqsmurU7%2Foh%2FG2kJZoNWJTiGsVPbnh48Wng6XwbEgLefI%3D&X-Amz-SignedHeaders=host&response-expires=Wed%2C%2004%20Jun%202025%2003%3A00%3A46%20GMT&X-Amz-Signature=ea1fa0e77bdd153cfc9eea96d390ae6538b81f2c9f525d109e58d2eaa936bb37'

$ws.Range("K3").Value = 'https://s3.samsara.com/samsara-dashcam-videos/4006124/281474991205821/1748885144009/VaUCDIqpi8-camera-video-segment-driver-1748885146509.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSB6HHK3ZB%2F20250603%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250603T190046Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEEMaCXVzLXdlc3QtMiJHMEUCIQDKu6SVHCuPxyLrQICvTWjsc8bJ1CoUDXThxrO8uE40HgIgA430WHEcncckntaqxq6OjxUxQq%2F4cRp0Tc2i6CEzSzcq3QMIGxAEGgw3ODEyMDQ5NDIyNDQiDI9ylrE874qsoNZPVyq6AyjSTgVQ3WieFhqFryrKxkcF3ppOLD9zcpUrjVHb3tfkBN6XzMndrmHyS91MRE81cdbF8AgqQTfi0m2KmkyBJ6n4LzAczOL5eQmIGQqjRTfVfv4MbIJNcX%2FmVKFDY46VPtANd2dKZKWtvOAG2jurMALhCaXaFkW38Wi%2FDpDX4kGeQYxsZ%2BWqA%2FsIP1TUnMfrtTnXOXr2FGA%2FMtA8a3pU%2B9b5iTKTFU1q8iQeuEehHKOXDFejJzr6rd39mMAWUoBhL7ZoAfvRYwq%2FK%2Bu5qH2XGR4aTPF6Dx2xIKkAwnhQPUa0qE2CeRnfYOMJhWzPQZ%2FCAldov9ZAih%2FeGjhBkztd%2BwikotUi51k1gb5djBDdD2%2BnrfS7%2FA%2BTSmkoMhyATJSRIqetQy1V%2B9I62PDv6sKYPAk1LKJnxA9P2VzZjctcIoYlmqoLPRDc2V8iEP0%2FlBt2b6DkkgBhIYGN%2Bc5rO4G%2FDic5glV%2BC219thdfzs8tXE%2BWfqMdD%2BC6aDJnVkR5UhAxhc7sAU4AkjlsvBnyBmxe132qO3Z%2FiJXr5kbp5dblgunXwgVgK9j%2B8yggr%2BiTKBbxM48Jy6PHKAR5CS4wyvr8wQY6pQGka52Q1n8aqPwGgfIPuwP%2B0wrpKmLdlRxGHPY4UeSOY6kBGwsVeZJGugNXCjDZUPIt8Kbou5blV%2Fh4c21o97PMLU%2FWg4Uln96yBi32j2p9lYTYOYzY6hmHQZ3%2Br1IP6tipUXxLrSP0lJlUGMnfhDJrO4lrGDrfWRWkkBBLnKeb%2FXTdLeeqsmurU7%2Foh%2FG2kJZoNWJTiGsVPbnh48Wng6XwbEgLefI%3D&X-Amz-SignedHeaders=host&response-expires=Wed%2C%2004%20Jun%202025%2003%3A00%3A46%20GMT&X-Amz-Signature=702634276eeab35e527cfb063249ef70ed9b5952627f3c3e1438c5b4731987fe'

$ws.Range("L3").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991205821/1748885144009/8zw63V5M8l-camera-video-segment-1748885146509.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSB6HHK3ZB%2F20250603%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250603T190046Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEEMaCXVzLXdlc3QtMiJHMEUCIQDKu6SVHCuPxyLrQICvTWjsc8bJ1CoUDXThxrO8uE40HgIgA430WHEcncckntaqxq6OjxUxQq%2F4cRp0Tc2i6CEzSzcq3QMIGxAEGgw3ODEyMDQ5NDIyNDQiDI9ylrE874qsoNZPVyq6AyjSTgVQ3WieFhqFryrKxkcF3ppOLD9zcpUrjVHb3tfkBN6XzMndrmHyS91MRE81cdbF8AgqQTfi0m2KmkyBJ6n4LzAczOL5eQmIGQqjRTfVfv4MbIJNcX%2FmVKFDY46VPtANd2dKZKWtvOAG2jurMALhCaXaFkW38Wi%2FDpDX4kGeQYxsZ%2BWqA%2FsIP1TUnMfrtTnXOXr2FGA%2FMtA8a3pU%2B9b5iTKTFU1q8iQeuEehHKOXDFejJzr6rd39mMAWUoBhL7ZoAfvRYwq%2FK%2Bu5qH2XGR4aTPF6Dx2xIKkAwnhQPUa0qE2CeRnfYOMJhWzPQZ%2FCAldov9ZAih%2FeGjhBkztd%2BwikotUi51k1gb5djBDdD2%2BnrfS7%2FA%2BTSmkoMhyATJSRIqetQy1V%2B9I62PDv6sKYPAk1LKJnxA9P2VzZjctcIoYlmqoLPRDc2V8iEP0%2FlBt2b6DkkgBhIYGN%2Bc5rO4G%2FDic5glV%2BC219thdfzs8tXE%2BWfqMdD%2BC6aDJnVkR5UhAxhc7sAU4AkjlsvBnyBmxe132qO3Z%2FiJXr5kbp5dblgunXwgVgK9j%2B8yggr%2BiTKBbxM48Jy6PHKAR5CS4wyvr8wQY6pQGka52Q1n8aqPwGgfIPuwP%2B0wrpKmLdlRxGHPY4UeSOY6kBGwsVeZJGugNXCjDZUPIt8Kbou5blV%2Fh4c21o97PMLU%2FWg4Uln96yBi32j2p9lYTYOYzY6hmHQZ3%2Br1IP6tipUXxLrSP0lJlUGMnfhDJrO4lrGDrfWRWkkBBLnKeb%2FXTdLeeqsmurU7%2Foh%2FG2kJZoNWJTiGsVPbnh48Wng6XwbEgLefI%3D&X-Amz-SignedHeaders=host&response-expires=Wed%2C%2004%20Jun%202025%2003%3A00%3A46%20GMT&X-Amz-Signature=bd8a47742b457adcd2c1286a6e62b04bc75c29b2af755ea56159e162b6abefe1'
